# Add a "skos:definition" column (E) to the AERO thesaurus sheet, with a
# header in E7 and the A320 family definition text in E10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the data cell before the header cell so the shared-string table ends
# up with the definition text at the lower index and the header label after
# it (matches the order the strings were originally authored in).
$ws.Range("E10").Value = "A320 Family Aircraft.  An A320 is a narrow body twin jet aircraft that is manufactured by Airbus"
$ws.Range("E7").Value = "skos:definition"

# E7 is a header cell like A7:D7, so give it the same bold header formatting.
$ws.Range("A7").Copy()
$ws.Range("E7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Both the header and the definition text should wrap within the column.
$ws.Range("E7").WrapText = $true
$ws.Range("E10").WrapText = $true

# Widen column E so the wrapped definition text is readable.
$ws.Columns("E").ColumnWidth = 35.5
